# Apply updated cryptocurrency data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.446.01'
$ws.Range('E2').Value = '  +2.74%  '

$ws.Range('D3').Value = '1.606.62'
$ws.Range('E3').Value = '  +2.71%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').Value = '''212.84'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('E6').Value = '  +7.01%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').Value = '''26.72'
$ws.Range('E8').Value = '  +6.44%  '

$ws.Range('D9').Value = '''43.61'
$ws.Range('E9').Value = '  -0.86%  '

$ws.Range('E10').Value = '  +2.75%  '

$ws.Range('E11').Value = '  +2.77%  '

$ws.Range('D12').Value = '''0.0911'
$ws.Range('E12').Value = '  +1.91%  '

$ws.Range('D13').Value = '1.837.62'
$ws.Range('E13').Value = '  +2.86%  '

$ws.Range('D14').Value = '1.617.48'
$ws.Range('E14').Value = '  +3.44%  '

$ws.Range('D15').Value = '29.468.15'
$ws.Range('E15').Value = '  +2.76%  '

$ws.Range('E16').Value = '  +3.63%  '

$ws.Range('E17').Value = '  +1.98%  '

$ws.Range('D18').Value = '''63.43'
$ws.Range('E18').Value = '  +3.50%  '

$ws.Range('D19').Value = '''240.30'
$ws.Range('E19').Value = '  +5.12%  '

$ws.Range('D20').Value = '''7.62'
$ws.Range('E20').Value = '  +3.78%  '

$ws.Range('E21').Value = '  +1.70%  '

$ws.Range('E22').Value = '  +0.25%  '

$ws.Range('D23').Value = '''4.00'
$ws.Range('E23').Value = '  +2.17%  '

$ws.Range('D24').Value = '''9.21'
$ws.Range('E24').Value = '  +2.05%  '

$ws.Range('E25').Value = '  +0.22%  '

$ws.Range('D26').Value = '''154.30'
$ws.Range('E26').Value = '  +1.97%  '

$ws.Range('E27').Value = '  +5.00%  '

$ws.Range('E28').Value = '  +3.34%  '

$ws.Range('E29').Value = '  +2.50%  '

$ws.Range('E30').Value = '  +0.20%  '

$ws.Range('E31').Value = '  +2.55%  '

$ws.Range('E32').Value = '  +1.44%  '

$ws.Range('D33').Value = '''3.23'
$ws.Range('E33').Value = '  +1.66%  '

$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.413.85'
$ws.Range('E34').Value = '  +1.47%  '

$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''3.09'
$ws.Range('E35').Value = '  +3.69%  '

$ws.Range('E36').Value = '  +0.94%  '

$ws.Range('E37').Value = '  +4.13%  '

$ws.Range('D38').Value = '''2.82'
$ws.Range('E38').Value = '  +4.63%  '

$ws.Range('E39').Value = '  +0.49%  '

$ws.Range('E40').Value = '  +2.64%  '

$ws.Range('D41').Value = '''0.537'
$ws.Range('E41').Value = '  +3.79%  '

$ws.Range('D42').Value = '''1.98'
$ws.Range('E42').Value = '  +1.66%  '

$ws.Range('D43').Value = '''0.0486'
$ws.Range('E43').Value = '  +5.99%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '''0.797'
$ws.Range('E44').Value = '  +3.39%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.21%  '

$ws.Range('D46').Value = '''52.73'
$ws.Range('E46').Value = '  +22.03%  '

$ws.Range('D47').Value = '''65.83'
$ws.Range('E47').Value = '  +2.88%  '

$ws.Range('E48').Value = '  +0.79%  '

$ws.Range('D49').Value = '1.746.75'
$ws.Range('E49').Value = '  +3.00%  '

$ws.Range('D50').Value = '''0.858'
$ws.Range('E50').Value = '  -1.49%  '

$ws.Range('D51').Value = '''86.62'
$ws.Range('E51').Value = '  +1.77%  '
